$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$s70 = @'
[{"ad":"<html><body><script type=\"text\/javascript\">var ut_ju = 'http://stgads.undertone.com/aj';ut = new Object();ut.bidid='21b46f0d859b33';ut.bannerid=1196449;ut.zoneid=137770;ut.hbprice=4.00;ut.width=1;ut.height=1;ut.adaptor='prebid';ut.pid=3470;ut.extpid='10433394';ut.domain='ranker.com';<\/script><script type=\"text/javascript\" src=\"http://cdn.undertone.com/js/ajs.js\"><\/script><\/body><\/html>","publisherId": 3470,"bidRequestId": "21b46f0d859b33","placementId": "10433394","adId": 1196449,"campaignId": 297398,"height": 1,"width": 1,"currency": "USD","cpm": 4,"adaptor": "prebid","netRevenue": "true"}]
'@

$s71 = @'
[{"ad":"<html><body><script type=\"text\/javascript\">var ut_ju = 'http://stgads.undertone.com/aj';ut = new Object();ut.bidid='21b46f0d859b99';ut.bannerid=1179961;ut.zoneid=152010;ut.hbprice=7.00;ut.width=1;ut.height=1;ut.adaptor='prebid';ut.pid=3470;ut.extpid='10433394';ut.domain='ranker.com';<\/script><script type=\"text/javascript\" src=\"http://cdn.undertone.com/js/ajs.js\"><\/script><\/body><\/html>","publisherId": 3470,"bidRequestId": "21b46f0d859b99","placementId": "10433394","adId": 1179961,"campaignId": 292801,"height": 1,"width": 1,"currency": "USD","cpm": 7,"adaptor": "prebid","netRevenue": "true"}]
'@

$s72 = @'
[{"ad":"<html><body><script type=\"text\/javascript\">var ut_ju = 'http://stgads.undertone.com/aj';ut = new Object();ut.bidid='21b46f0d859b33';ut.bannerid=1196449;ut.zoneid=137770;ut.hbprice=4.00;ut.width=1;ut.height=1;ut.adaptor='prebid';ut.pid=3470;ut.extpid='';ut.domain='ranker.com';<\/script><script type=\"text/javascript\" src=\"http://cdn.undertone.com/js/ajs.js\"><\/script><\/body><\/html>","publisherId": 3470,"bidRequestId": "21b46f0d859b33","placementId": "","adId": 1196449,"campaignId": 297398,"height": 1,"width": 1,"currency": "USD","cpm": 4,"adaptor": "prebid","netRevenue": "true"}]
'@

$s73 = @'
[{"ad":"<html><body><script type=\"text\/javascript\">var ut_ju = 'http://stgads.undertone.com/aj';ut = new Object();ut.bidid='21b46f0d859b47';ut.bannerid=1196449;ut.zoneid=137770;ut.hbprice=4.00;ut.width=1;ut.height=1;ut.adaptor='prebid';ut.pid=3470;ut.extpid='10433394';ut.domain='ranker.com';<\/script><script type=\"text/javascript\" src=\"http://cdn.undertone.com/js/ajs.js\"><\/script><\/body><\/html>","publisherId": 3470,"bidRequestId": "21b46f0d859b47","placementId": "10433394","adId": 1196449,"campaignId": 297398,"height": 1,"width": 1,"currency": "USD","cpm": 4,"adaptor": "prebid","netRevenue": "true"}]
'@

# Update Expected results (column F) for the ranker.com test rows
$ws.Range("F13").Value = $s70
$ws.Range("F15").Value = $s70
$ws.Range("F16").Value = $s71
$ws.Range("F18").Value = $s72
$ws.Range("F22").Value = $s73

# Match the author workbook view state: scrolled down with F22 selected
$ws.Range("F22").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
